# Update performance dashboard 2025-12-24 23:49 - Simplified Design
#
# For each sheet:
#  - Insert a new column N ("Files Count") before the existing "Analysis Date"
#    column, pushing "Analysis Date" to O and "Data Date" to P.
#  - Populate the new "Files Count" header + per-row values.
#  - Update "Trading Days" (M) and "Annual Return (%)" (G) values to reflect
#    the re-run analysis.

$wb = $excel.ActiveWorkbook

# row layout (applies to every sheet): row, AnnualReturn(G), TradingDays/FilesCount(M/N)
$rowsSummary = @(
    @{ Row = 2;  G = "+216.99%"; Days = 6 },
    @{ Row = 3;  G = "+195.90%"; Days = 6 },
    @{ Row = 4;  G = "+51.47%";  Days = 6 },
    @{ Row = 5;  G = "+33.17%";  Days = 6 },
    @{ Row = 6;  G = "+3.27%";   Days = 6 },
    @{ Row = 7;  G = "+17.36%";  Days = 6 },
    @{ Row = 8;  G = "+108.09%"; Days = 6 },
    @{ Row = 9;  G = "+17.68%";  Days = 6 },
    @{ Row = 10; G = "+5.95%";   Days = 6 },
    @{ Row = 11; G = "+2.14%";   Days = 6 },
    @{ Row = 12; G = "+16.98%";  Days = 5 },
    @{ Row = 13; G = "+134.88%"; Days = 6 },
    @{ Row = 14; G = "+126.49%"; Days = 6 },
    @{ Row = 15; G = "+3.66%";   Days = 6 },
    @{ Row = 16; G = "+39.58%";  Days = 6 }
)

$rowsPattern1 = @(
    @{ Row = 2; G = "+216.99%"; Days = 6 },
    @{ Row = 3; G = "+195.90%"; Days = 6 },
    @{ Row = 4; G = "+51.47%";  Days = 6 },
    @{ Row = 5; G = "+33.17%";  Days = 6 },
    @{ Row = 6; G = "+3.27%";   Days = 6 }
)

$rowsPattern2 = @(
    @{ Row = 2; G = "+17.36%";  Days = 6 },
    @{ Row = 3; G = "+108.09%"; Days = 6 },
    @{ Row = 4; G = "+17.68%";  Days = 6 },
    @{ Row = 5; G = "+5.95%";   Days = 6 },
    @{ Row = 6; G = "+2.14%";   Days = 6 }
)

$rowsPattern3 = @(
    @{ Row = 2; G = "+16.98%";  Days = 5 },
    @{ Row = 3; G = "+134.88%"; Days = 6 },
    @{ Row = 4; G = "+126.49%"; Days = 6 },
    @{ Row = 5; G = "+3.66%";   Days = 6 },
    @{ Row = 6; G = "+39.58%";  Days = 6 }
)

$sheetPlans = @(
    @{ Name = "Summary";                   Rows = $rowsSummary },
    @{ Name = "Pattern1-Pure Data";        Rows = $rowsPattern1 },
    @{ Name = "Pattern2-Data+Technical";   Rows = $rowsPattern2 },
    @{ Name = "Pattern3-Data+News";        Rows = $rowsPattern3 }
)

foreach ($plan in $sheetPlans) {
    $ws = $wb.Worksheets.Item($plan.Name)

    # Insert new column at N (14): shifts old N (Analysis Date) -> O,
    # old O (Data Date) -> P.
    $ws.Columns.Item(14).Insert()

    # New header cell for the inserted column.
    $ws.Cells.Item(1, 14).Value2 = "Files Count"

    foreach ($r in $plan.Rows) {
        $rowNum = $r.Row

        # Annual Return (%) - column G (7). This column holds a literal
        # text percentage (e.g. "+216.99%"), not a real Excel percentage
        # number, so force text format before assigning then drop the
        # number-format override again so the cell keeps the plain
        # (unstyled) look the other literal-text columns use.
        $gCell = $ws.Cells.Item($rowNum, 7)
        $gCell.NumberFormat = "@"
        $gCell.Value2 = $r.G
        $gCell.ClearFormats()

        # Trading Days - column M (13)
        $ws.Cells.Item($rowNum, 13).Value2 = $r.Days
        # Files Count (new) - column N (14)
        $ws.Cells.Item($rowNum, 14).Value2 = $r.Days
    }
}
